$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18: new timesheet entry (restart of project)
$ws.Range("A18").Value = "2016-08-10"
$ws.Range("D18").Value = "<- EK2 restart"
$ws.Range("B18").Value = "opnieuw beginnen: databases, views, controllers"
$ws.Range("C18").Value = 8

# Row 19: date only
$ws.Range("A19").Value = "2016-08-11"

# Update selected cell to match the saved view state
[void]$ws.Range("B19").Select()
